# Fruta / hortaliza, semanal
# Insert a new weekly record at row 52 ("Feria Lagunitas de Puerto Montt - Poroto verde"),
# pushing all subsequent rows (old 52..84) down by one (new 53..85).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 52 and below down by one row.
$ws.Rows(52).Insert()

# Populate the newly inserted row 52 with the new weekly observation.
$ws.Range("A52").Value = 4
$ws.Range("B52").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C52").Value = "Los Lagos"
$ws.Range("D52").Value = 44729
$ws.Range("E52").Value = 10
$ws.Range("F52").Value = 100112031
$ws.Range("G52").Value = "Poroto verde"
$ws.Range("H52").Value = "Magnum"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 40
$ws.Range("K52").Value = 26000
$ws.Range("L52").Value = 26000
$ws.Range("M52").Value = 26000
$ws.Range("N52").Value = "$/malla 25 kilos"
$ws.Range("O52").Value = "Perú"
$ws.Range("P52").Value = 1040
$ws.Range("Q52").Value = 25
$ws.Range("R52").Value = "Hortaliza"
